$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("Главные", "Линейные")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("AA2:AA26").Value = "2025-12-14 07:02:27"
}
